# Auto-generated update of market/profit lookup values (H,I,J,K,L,M,N columns)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, per scheduled data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1154.375
$ws.Range("I28").Value = 564.73334
$ws.Range("K28").Value = 564.73334
$ws.Range("M28").Value = -79.73334
$ws.Range("H40").Value = 1995.625
$ws.Range("I40").Value = 1995.7142
$ws.Range("J40").Value = 1995
$ws.Range("K40").Value = 1995.7142
$ws.Range("L40").Value = 1995
$ws.Range("M40").Value = -1820.7142
$ws.Range("N40").Value = -2345
$ws.Range("H86").Value = 2693.3635
$ws.Range("J86").Value = 3004.6
$ws.Range("L86").Value = 3004.6
$ws.Range("N86").Value = -5250.6
$ws.Range("H89").Value = 2693.3635
$ws.Range("J89").Value = 3004.6
$ws.Range("L89").Value = 15023
$ws.Range("N89").Value = -26255
$ws.Range("H106").Value = 2749.8333
$ws.Range("I106").Value = 2499.8
$ws.Range("K106").Value = 2499.8
$ws.Range("M106").Value = -1868.8
$ws.Range("H113").Value = 6106
$ws.Range("I113").Value = 7329.4287
$ws.Range("J113").Value = 3251.3333
$ws.Range("K113").Value = 7329.4287
$ws.Range("L113").Value = 3251.3333
$ws.Range("M113").Value = -4075.4287
$ws.Range("N113").Value = -9759.3333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 21700.033
$ws.Range("J32").Value = 8999
$ws.Range("L32").Value = 8999
$ws.Range("N32").Value = -9573
$ws.Range("H74").Value = 334492.78
$ws.Range("I74").Value = 401111.47
$ws.Range("K74").Value = 401111.47
$ws.Range("M74").Value = -400237.47
$ws.Range("H77").Value = 334492.78
$ws.Range("I77").Value = 401111.47
$ws.Range("K77").Value = 2005557.35
$ws.Range("M77").Value = -2001189.35
$ws.Range("H122").Value = 1515.5
$ws.Range("I122").Value = 1194.375
$ws.Range("J122").Value = 2800
$ws.Range("K122").Value = 3583.125
$ws.Range("L122").Value = 8400
$ws.Range("M122").Value = -1133.125
$ws.Range("N122").Value = -13300
$ws.Range("H132").Value = 2898.1667
$ws.Range("I132").Value = 1748.8
$ws.Range("J132").Value = 4813.778
$ws.Range("K132").Value = 5246.4
$ws.Range("L132").Value = 14441.334
$ws.Range("M132").Value = -2716.4
$ws.Range("N132").Value = -19501.334
$ws.Range("H140").Value = 59664.4
$ws.Range("J140").Value = 59664.4
$ws.Range("L140").Value = 59664.4
$ws.Range("N140").Value = -70024.39999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 16833.281
$ws.Range("J20").Value = 1148.7778
$ws.Range("L20").Value = 1148.7778
$ws.Range("N20").Value = -1642.7778
$ws.Range("H74").Value = 81441.5
$ws.Range("J74").Value = 81441.5
$ws.Range("L74").Value = 81441.5
$ws.Range("N74").Value = -83313.5
$ws.Range("H77").Value = 81441.5
$ws.Range("J77").Value = 81441.5
$ws.Range("L77").Value = 244324.5
$ws.Range("N77").Value = -253684.5
$ws.Range("H94").Value = 1426.3462
$ws.Range("I94").Value = 798.1579
$ws.Range("J94").Value = 3131.4285
$ws.Range("K94").Value = 798.1579
$ws.Range("L94").Value = 3131.4285
$ws.Range("M94").Value = -347.1579
$ws.Range("N94").Value = -4033.4285
$ws.Range("H105").Value = 5668.143
$ws.Range("I105").Value = 4084
$ws.Range("J105").Value = 8242.375
$ws.Range("K105").Value = 4084
$ws.Range("L105").Value = 8242.375
$ws.Range("M105").Value = -2337
$ws.Range("N105").Value = -11736.375
$ws.Range("H107").Value = 29640.611
$ws.Range("J107").Value = 1387.8
$ws.Range("L107").Value = 1387.8
$ws.Range("N107").Value = -5227.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H133").Value = 116324.93
$ws.Range("J133").Value = 115714
$ws.Range("L133").Value = 115714
$ws.Range("N133").Value = -120774
$ws.Range("H137").Value = 87121.2
$ws.Range("J137").Value = 87121.2
$ws.Range("L137").Value = 87121.2
$ws.Range("N137").Value = -97321.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1371.5
$ws.Range("I113").Value = 643.6667
$ws.Range("J113").Value = 1735.4166
$ws.Range("K113").Value = 1931.0001
$ws.Range("L113").Value = 5206.2498
$ws.Range("M113").Value = 238.9999
$ws.Range("N113").Value = -9546.2498
$ws.Range("H132").Value = 2197.9167
$ws.Range("I132").Value = 2083.7778
$ws.Range("J132").Value = 2266.4
$ws.Range("K132").Value = 18754.0002
$ws.Range("L132").Value = 20397.6
$ws.Range("M132").Value = -16224.0002
$ws.Range("N132").Value = -25457.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6443.5
$ws.Range("I70").Value = 6225.5415
$ws.Range("K70").Value = 6225.5415
$ws.Range("M70").Value = -5955.5415
$ws.Range("H73").Value = 6443.5
$ws.Range("I73").Value = 6225.5415
$ws.Range("K73").Value = 6225.5415
$ws.Range("M73").Value = -5289.5415
$ws.Range("H97").Value = 879.98
$ws.Range("I97").Value = 703.8
$ws.Range("J97").Value = 1584.7
$ws.Range("K97").Value = 703.8
$ws.Range("L97").Value = 1584.7
$ws.Range("M97").Value = -207.8
$ws.Range("N97").Value = -2576.7
$ws.Range("H140").Value = 127966.336
$ws.Range("J140").Value = 127966.336
$ws.Range("L140").Value = 127966.336
$ws.Range("N140").Value = -138326.336

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3131.9092
$ws.Range("I7").Value = 2775.125
$ws.Range("K7").Value = 2775.125
$ws.Range("M7").Value = -2663.125
$ws.Range("H16").Value = 2923.625
$ws.Range("I16").Value = 2948.3333
$ws.Range("K16").Value = 2948.3333
$ws.Range("M16").Value = -2778.3333
$ws.Range("H46").Value = 5368.5557
$ws.Range("I46").Value = 975.1667
$ws.Range("J46").Value = 7565.25
$ws.Range("K46").Value = 975.1667
$ws.Range("L46").Value = 7565.25
$ws.Range("M46").Value = -787.1667
$ws.Range("N46").Value = -7941.25
$ws.Range("H70").Value = 26163
$ws.Range("J70").Value = 26163
$ws.Range("L70").Value = 26163
$ws.Range("N70").Value = -26703
$ws.Range("H73").Value = 26163
$ws.Range("J73").Value = 26163
$ws.Range("L73").Value = 26163
$ws.Range("N73").Value = -28035
$ws.Range("H82").Value = 10411.533
$ws.Range("I82").Value = 9523.299999999999
$ws.Range("J82").Value = 12188
$ws.Range("K82").Value = 9523.299999999999
$ws.Range("L82").Value = 12188
$ws.Range("M82").Value = -9162.299999999999
$ws.Range("N82").Value = -12910
$ws.Range("H85").Value = 10411.533
$ws.Range("I85").Value = 9523.299999999999
$ws.Range("J85").Value = 12188
$ws.Range("K85").Value = 9523.299999999999
$ws.Range("L85").Value = 12188
$ws.Range("M85").Value = -8275.299999999999
$ws.Range("N85").Value = -14684
$ws.Range("H96").Value = 70000
$ws.Range("J96").Value = 70000
$ws.Range("L96").Value = 70000
$ws.Range("N96").Value = -75492
$ws.Range("H126").Value = 3131.9092
$ws.Range("I126").Value = 2775.125
$ws.Range("K126").Value = 8325.375
$ws.Range("M126").Value = -5855.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 14359.4
$ws.Range("I47").Value = 14999.25
$ws.Range("J47").Value = 11800
$ws.Range("K47").Value = 14999.25
$ws.Range("L47").Value = 11800
$ws.Range("M47").Value = -14427.25
$ws.Range("N47").Value = -12944
$ws.Range("H86").Value = 51666.668
$ws.Range("J86").Value = 51666.668
$ws.Range("L86").Value = 51666.668
$ws.Range("N86").Value = -53912.668
$ws.Range("H89").Value = 51666.668
$ws.Range("J89").Value = 51666.668
$ws.Range("L89").Value = 258333.34
$ws.Range("N89").Value = -269565.34
$ws.Range("H99").Value = 64833.332
$ws.Range("J99").Value = 64833.332
$ws.Range("L99").Value = 64833.332
$ws.Range("N99").Value = -70823.33199999999
$ws.Range("H120").Value = 50420
$ws.Range("J120").Value = 50420
$ws.Range("L120").Value = 50420
$ws.Range("N120").Value = -60096
$ws.Range("H126").Value = 627637.6
$ws.Range("I126").Value = 2683.5
$ws.Range("J126").Value = 2502500
$ws.Range("K126").Value = 8050.5
$ws.Range("L126").Value = 7507500
$ws.Range("M126").Value = -5580.5
$ws.Range("N126").Value = -7512440
$ws.Range("H132").Value = 2619.7886
$ws.Range("I132").Value = 2614.4546
$ws.Range("J132").Value = 2649.125
$ws.Range("K132").Value = 7843.3638
$ws.Range("L132").Value = 7947.375
$ws.Range("M132").Value = -5313.3638
$ws.Range("N132").Value = -13007.375
$ws.Range("H136").Value = 31183.959
$ws.Range("I136").Value = 37636.367
$ws.Range("J136").Value = 6664.8
$ws.Range("K136").Value = 112909.101
$ws.Range("L136").Value = 19994.4
$ws.Range("M136").Value = -110359.101
$ws.Range("N136").Value = -25094.4

